$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (they contain dotted/decimal
# strings that must not be re-interpreted as numbers by Excel).
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D19","D20","D21","D23","D24","D25","D26","D28","D30","D31","D32","D33","D35","D36","D39","D40","D41","D42","D43","D45","D46","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped from the source diff.
$ws.Range("D2").Value = "30.394.10"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.876.01"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "238.71"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.4792"
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("D8").Value = "0.2819"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").Value = "0.06514"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").Value = "1.874.28"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "0.07478"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "16.62"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "5.073"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "88.20"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "0.6599"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "30.369.73"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "0.000007578"
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").Value = "2.116.73"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").Value = "5.300"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "220.25"
$ws.Range("E23").Value = "  +13.86%  "
$ws.Range("D24").Value = "6.187"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "9.352"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "167.77"
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").Value = "1.966"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").Value = "0.09352"
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("D31").Value = "4.312"
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").Value = "4.025"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "0.05026"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("E34").Value = "  +4.08%  "
$ws.Range("D35").Value = "0.7436"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "2.709"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "0.9056"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "2.063"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "106.53"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "5.879"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "0.4274"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "7.406"
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("D46").Value = "64.64"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("E47").Value = "  -4.08%  "
$ws.Range("D48").Value = "1.472"
$ws.Range("E48").Value = "  -5.71%  "
$ws.Range("D49").Value = "8.921"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "33.75"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "0.3885"
$ws.Range("E51").Value = "  +0.51%  "
